# Apply "new sensitivity and calculus" changes to the daily model results workbook.
# Sheet1 "Model Accuracy": add new columns (Market threshold, Market min, Market max,
#   Recall, Precision) and refresh Accuracy values.
# Sheets2-6 "Confusion Matrix ...": refresh a handful of cell counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Model Accuracy (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Copy the formatting of the existing bold/bordered header (B1) into the new
# header cells, then overwrite the copied value with the correct text.
$ws1.Range("B1").Copy($ws1.Range("C1"))
$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("B1").Copy($ws1.Range("D1"))
$ws1.Range("D1").Value = "Market min"
$ws1.Range("B1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "Market max"
$ws1.Range("B1").Copy($ws1.Range("F1"))
$ws1.Range("F1").Value = "Recall"
$ws1.Range("B1").Copy($ws1.Range("G1"))
$ws1.Range("G1").Value = "Precision"

# Row 2 - TOTALENERGIES SE
$ws1.Range("B2").Value = 42.6039119804401
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 11.11111111111111
$ws1.Range("G2").Value = 1.333333333333333

# Row 3 - FMC CORP
$ws1.Range("B3").Value = 28.54523227383863
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 12.8686327077748
$ws1.Range("G3").Value = 25.66844919786097

# Row 4 - BP PLC
$ws1.Range("B4").Value = 65.89242053789731
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

# Row 5 - STORA ENSO
$ws1.Range("B5").Value = 52.50611246943765
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 7.272727272727272
$ws1.Range("G5").Value = 4.761904761904762

# Row 6 - BHP GROUP
$ws1.Range("B6").Value = 61.06356968215159
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 50
$ws1.Range("G6").Value = 0.5434782608695652

# ---------------------------------------------------------------------------
# Sheet 2: Confusion Matrix TOTALENERGIES SE (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C3").Value = 689
$ws2.Range("D3").Value = 3
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 336

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix FMC CORP (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 48
$ws3.Range("C2").Value = 98
$ws3.Range("D2").Value = 41
$ws3.Range("B3").Value = 146
$ws3.Range("C3").Value = 259
$ws3.Range("D3").Value = 155
$ws3.Range("B4").Value = 179
$ws3.Range("C4").Value = 295
$ws3.Range("D4").Value = 160

# ---------------------------------------------------------------------------
# Sheet 4: Confusion Matrix BP PLC (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B3").Value = 25
$ws4.Range("C3").Value = 1063
$ws4.Range("D3").Value = 27
$ws4.Range("B4").Value = 15
$ws4.Range("C4").Value = 424
$ws4.Range("D4").Value = 15

# ---------------------------------------------------------------------------
# Sheet 5: Confusion Matrix STORA ENSO (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = 8
$ws5.Range("C2").Value = 152
$ws5.Range("D2").Value = 8
$ws5.Range("B3").Value = 70
$ws5.Range("C3").Value = 819
$ws5.Range("D3").Value = 67
$ws5.Range("B4").Value = 32
$ws5.Range("C4").Value = 387
$ws5.Range("D4").Value = 32

# ---------------------------------------------------------------------------
# Sheet 6: Confusion Matrix BHP GROUP (-0.15, 0.15, 0.15)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B2").Value = 2
$ws6.Range("C2").Value = 365
$ws6.Range("D2").Value = 1
$ws6.Range("B3").Value = 2
$ws6.Range("C3").Value = 997
$ws6.Range("D3").Value = 2
$ws6.Range("B4").Value = 0
$ws6.Range("C4").Value = 211
$ws6.Range("D4").Value = 0
